$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "Monday, Jan 09"
$ws.Range("C51").Value = "1:55 PM"
$ws.Range("D51").Value = "SK1755"
$ws.Range("E51").Value = "Copenhagen"
$ws.Range("F51").Value = "(CPH)"
$ws.Range("G51").Value = "SAS "
$ws.Range("H51").Value = "CRJ9"
$ws.Range("I51").Value = "(EI-FPV)"
$ws.Range("J51").Value = "1:55 PM"
$ws.Range("L51").Value = "0 hours, 0 minutes"
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "Monday, Jan 09"
$ws.Range("C52").Value = "2:15 PM"
$ws.Range("D52").Value = "LO3943"
$ws.Range("E52").Value = "Warsaw"
$ws.Range("F52").Value = "(WAW)"
$ws.Range("G52").Value = "LOT "
$ws.Range("H52").Value = "E170"
$ws.Range("I52").Value = "(SP-LDG)"
$ws.Range("J52").Value = "2:02 PM"
$ws.Range("L52").Value = "0 hours, -13 minutes"
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "Monday, Jan 09"
$ws.Range("C53").Value = "2:24 PM"
$ws.Range("D53").Value = "UNKNOWN"
$ws.Range("E53").Value = "Szczecin"
$ws.Range("F53").Value = "(SZZ)"
$ws.Range("G53").Value = "AMC Aviation "
$ws.Range("H53").Value = "PC24"
$ws.Range("I53").Value = "(SP-AGA)"
$ws.Range("J53").Value = "1:50 PM"
$ws.Range("L53").Value = "0 hours, -34 minutes"
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "Monday, Jan 09"
$ws.Range("C54").Value = "2:45 PM"
$ws.Range("D54").Value = "KL1273"
$ws.Range("E54").Value = "Amsterdam"
$ws.Range("F54").Value = "(AMS)"
$ws.Range("G54").Value = "KLM "
$ws.Range("H54").Value = "E75L"
$ws.Range("I54").Value = "(PH-EXW)"
$ws.Range("J54").Value = "2:30 PM"
$ws.Range("L54").Value = "0 hours, -15 minutes"
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "Monday, Jan 09"
$ws.Range("C55").Value = "3:30 PM"
$ws.Range("D55").Value = "FR7947"
$ws.Range("E55").Value = "Oslo"
$ws.Range("F55").Value = "(TRF)"
$ws.Range("G55").Value = "Ryanair "
$ws.Range("H55").Value = "B738"
$ws.Range("I55").Value = "(SP-RSX)"
$ws.Range("J55").Value = "3:17 PM"
$ws.Range("L55").Value = "0 hours, -13 minutes"
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "Monday, Jan 09"
$ws.Range("C56").Value = "4:00 PM"
$ws.Range("D56").Value = "W91902"
$ws.Range("E56").Value = "London"
$ws.Range("F56").Value = "(LTN)"
$ws.Range("G56").Value = "Wizz Air "
$ws.Range("H56").Value = "A320"
$ws.Range("I56").Value = "(G-WUKD)"
$ws.Range("J56").Value = "3:46 PM"
$ws.Range("L56").Value = "0 hours, -14 minutes"
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "Monday, Jan 09"
$ws.Range("C57").Value = "5:20 PM"
$ws.Range("D57").Value = "FR7890"
$ws.Range("E57").Value = "Malta"
$ws.Range("F57").Value = "(MLA)"
$ws.Range("G57").Value = "Ryanair "
$ws.Range("H57").Value = "B738"
$ws.Range("I57").Value = "(SP-RKR)"
$ws.Range("J57").Value = "5:29 PM"
$ws.Range("L57").Value = "0 hours, 9 minutes"
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "Monday, Jan 09"
$ws.Range("C58").Value = "5:40 PM"
$ws.Range("D58").Value = "FR7943"
$ws.Range("E58").Value = "Manchester"
$ws.Range("F58").Value = "(MAN)"
$ws.Range("G58").Value = "Ryanair "
$ws.Range("H58").Value = "B738"
$ws.Range("I58").Value = "(SP-RSM)"
$ws.Range("J58").Value = "5:25 PM"
$ws.Range("L58").Value = "0 hours, -15 minutes"
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "Monday, Jan 09"
$ws.Range("C59").Value = "5:50 PM"
$ws.Range("D59").Value = "LO3947"
$ws.Range("E59").Value = "Warsaw"
$ws.Range("F59").Value = "(WAW)"
$ws.Range("G59").Value = "LOT "
$ws.Range("H59").Value = "E170"
$ws.Range("I59").Value = "(SP-LDH)"
$ws.Range("J59").Value = "5:34 PM"
$ws.Range("L59").Value = "0 hours, -16 minutes"
